# Update countries & provincias Spain
# Applies the 5-Aug-2020 14:30 COVID data refresh to the "Pais" sheet:
#  - updates the "Datos actualizados..." timestamp
#  - refreshes case numbers for a number of countries
#  - a few countries swap ranking order (their row now shows a different
#    country name because the updated totals re-sorted them)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 14:30"

# --- Helper: write a full data row (country name + B:H numbers) ---
function Set-CountryRow {
    param(
        [int]$Row,
        [string]$Country,
        [double]$CasosTotales,
        [double]$NuevosCasos,
        [double]$CasosActivos,
        [double]$Recuperados,
        [double]$CasosCriticos,
        [double]$MuertesHoy,
        [double]$Muertes
    )
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Row 4 - Estados Unidos (values refreshed, same country)
Set-CountryRow 4 "Estados Unidos" 4919116 696 2482899 2275882 0 45 160335

# Row 6 - India (values refreshed, same country)
Set-CountryRow 6 "India" 1910795 4182 1282929 588010 0 36 39856

# Rows 40/41 - Kuwait moves ahead of Panama
Set-CountryRow 40 "Kuwait" 69425 651 60906 8051 0 3 468
Set-CountryRow 41 "Panama" 69424 0 43330 24572 0 0 1522

# Row 45 - Paises Bajos (values refreshed, same country)
Set-CountryRow 45 "Paises Bajos" 56381 426 0 0 0 3 6153

# Rows 55/56 - Ghana moves ahead of Kirguistan
Set-CountryRow 55 "Ghana" 39075 1263 35563 3313 0 8 199
Set-CountryRow 56 "Kirguistan" 38110 569 29513 7159 0 11 1438

# Row 59 - Azerbaiyan (values refreshed, same country)
Set-CountryRow 59 "Azerbaiyan" 33103 193 28840 3787 0 3 476

# Row 69 - Nepal (values refreshed, same country)
Set-CountryRow 69 "Nepal" 21390 381 15156 6174 0 2 60

# Row 78 - Dinamarca (values refreshed, same country)
Set-CountryRow 78 "Dinamarca" 14185 112 12753 816 0 0 616

# Rows 79/80 - Bosnia y Herzegovina moves ahead of Estado de Palestina
Set-CountryRow 79 "Bosnia y Herzegovina" 13138 282 6839 5920 0 6 379
Set-CountryRow 80 "Estado de Palestina" 13065 295 6618 6359 0 2 88

# Row 82 - Madagascar (values refreshed, same country)
Set-CountryRow 82 "Madagascar" 12222 327 9798 2297 0 4 127

# Row 87 - Consejo Danes para los Refugiados (values refreshed, same country)
Set-CountryRow 87 "Consejo Danes para los Refugiados" 9253 75 7821 1217 0 0 215

# Row 93 - Finlandia (values refreshed, same country)
Set-CountryRow 93 "Finlandia" 7512 29 6980 201 0 0 331

# Row 100 - Croacia (values refreshed, same country)
Set-CountryRow 100 "Croacia" 5376 58 4589 633 0 0 154

# Row 115 - Tailandia (values refreshed, same country)
Set-CountryRow 115 "Tailandia" 3328 7 3144 126 0 0 58

# Rows 159/160/161 - Vietnam moves ahead of Crucero and San Marino
Set-CountryRow 159 "Vietnam" 713 41 381 324 0 0 8
Set-CountryRow 160 "Crucero" 712 0 651 48 0 0 13
Set-CountryRow 161 "San Marino" 699 0 657 0 0 0 42
